$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("January-2021")

# New retailer order entries recorded for 10/11 Jan 2021 (columns Q/R)
$ws.Range("Q5").Value = 1040
$ws.Range("Q7").Value = 2080
$ws.Range("R7").Value = 2080
$ws.Range("R11").Value = 2080
$ws.Range("Q14").Value = 5200
$ws.Range("R21").Value = 5200
$ws.Range("Q31").Value = 2080
$ws.Range("R35").Value = 3120
$ws.Range("R46").Value = 3120
$ws.Range("R57").Value = 5200
$ws.Range("R59").Value = 2080
$ws.Range("Q61").Value = 2080
$ws.Range("R62").Value = 2080
$ws.Range("Q65").Value = 5200
$ws.Range("Q66").Value = 3120
$ws.Range("Q67").Value = 2080

# Reflect the author's scroll position / active cell at save time
$ws.Range("D44").Select()
